# "Multa fija agua externa"
#
# On the "Escalares" sheet, parameter P (row 7) used to represent the
# cost of external water, described as "Costo agua externa (US$/m³)".
# This change splits that into two parameters:
#   - P   : "Costo de bombear agua externa (US$/m³)" (keeps the 5.45 value)
#   - Pf  : "Costo fijo de uso de agua externa (US$)" (new row, value 200000)
#
# The new Pf row is inserted right after P (pushing Pmax..Mbig down by one
# row), and the sheet's selection ends up on D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escalares")

# Insert a new row 8, shifting the existing rows 8-13 down to 9-14.
$ws.Rows("8:8").Insert()

# Bring over the formatting (styles) from row 7 into the newly created row 8.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)

# Row 8 (new): parameter "Pf" - fixed cost of using external water.
$ws.Range("A8").Value2 = "Pf"

# Row 7: relabel P's description (it is now the pumping cost, not the
# previous combined external-water cost description).
$ws.Range("B7").Value2 = "Costo de bombear agua externa (US`$/m³)"

$ws.Range("B8").Value2 = "Costo fijo de uso de agua externa (US`$)"
$ws.Range("C8").Value2 = 200000

$ws.Range("D10").Select()
